$d = $word.ActiveDocument

# The document currently has a single (Normal-style) paragraph:
# "Test one two three". Grow it into the target 9-paragraph chain:
#   1 (empty) / 2 Heading One / 3 Paragraph one / 4 (empty) /
#   5 Paragraph Two / 6 (empty) / 7 Heading Two / 8 Para three / 9 Para four
#
# Approach:
#  - First insert a blank paragraph *before* the existing paragraph, so
#    the original paragraph's text can simply be overwritten (rather than
#    ever assigning an explicit "" to a paragraph -- leaving paragraphs
#    untouched keeps them free of an explicit <w:t/> element).
#  - Then grow the chain downward with InsertParagraphAfter, setting text
#    on each new paragraph as it's created. Blank paragraphs are left
#    completely untouched (no Range.Text assignment at all).
#  - New paragraphs inherit the style of the paragraph they were split
#    from, so all splitting happens first while everything is still
#    plain/Normal, and Heading1 is applied to the two heading paragraphs
#    only afterwards -- that way the heading style never "leaks" into
#    paragraphs split off after it.

$orig = $d.Paragraphs.Item(1)
$orig.Range.InsertParagraphBefore()      # 1: empty paragraph

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Heading One"           # 2: Heading One
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "Paragraph one"         # 3: Paragraph one
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4)              # 4: empty paragraph (left blank)
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "Paragraph Two"         # 5: Paragraph Two
$p5.Range.InsertParagraphAfter()

$p6 = $d.Paragraphs.Item(6)              # 6: empty paragraph (left blank)
$p6.Range.InsertParagraphAfter()

$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "Heading Two"           # 7: Heading Two
$p7.Range.InsertParagraphAfter()

$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Para three"            # 8: Para three
$p8.Range.InsertParagraphAfter()

$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Para four"             # 9: Para four

# Now apply the Heading1 style to the two heading paragraphs only.
$p2.Style = "Heading1"
$p7.Style = "Heading1"

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
